$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update tracking dates for row 7 (item 2) to 25-03-2015 (serial 42088)
$ws.Range("C7").Value = 42088
$ws.Range("D7").Value = 42088
$ws.Range("E7").Value = 42088

# Mark observation with "unidad 2" status text
$ws.Range("F7").Value = "Ya puede ser revisado por María Clemencia"

# Row height grows to accommodate wrapped text
$ws.Rows.Item(7).RowHeight = 29.25

# Selection moves to G7
$ws.Range("G7").Select()
